# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear stray/buggy cells that should no longer hold a value
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

# Corrected YoY forecast values (recomputed after fixing the naive component forecaster bug)
$ws.Range("E3").Value  = 3.386383090739975
$ws.Range("C4").Value  = 2.533533936850585
$ws.Range("E4").Value  = 0.9842934829757288
$ws.Range("E5").Value  = 3.612753212925446
$ws.Range("C6").Value  = 1.21254482274098
$ws.Range("E6").Value  = 2.158838189283219
$ws.Range("E7").Value  = 1.194058515117313
$ws.Range("C8").Value  = 0.4712609263772816
$ws.Range("E8").Value  = 1.409662779709797
$ws.Range("E10").Value = 4.595879021798344
$ws.Range("C11").Value = 4.109890522944326
$ws.Range("E12").Value = 0.02883756256673031
$ws.Range("E13").Value = 0.9262553939923146
$ws.Range("E14").Value = 2.928189816005689
$ws.Range("E15").Value = 3.828814763561783
$ws.Range("C16").Value = 2.777797690741446
$ws.Range("C18").Value = -1.432689847121826
$ws.Range("E18").Value = 0.1752798163574321
$ws.Range("C19").Value = 2.033479419175155

$wb.Save()
